# Update the date and the two-digit multiplication problems to the new
# values generated for the updated output (commit 1c8df47).

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2023-08-13 Sunday" "2023-08-14 Monday"

Replace-Text "11×34=" "79×31="
Replace-Text "58×90=" "69×51="
Replace-Text "91×78=" "96×57="
Replace-Text "70×50=" "79×96="
Replace-Text "27×82=" "95×23="

Replace-Text "87×49=" "37×16="
Replace-Text "40×57=" "48×21="
Replace-Text "63×37=" "29×63="
Replace-Text "52×27=" "65×54="
Replace-Text "74×35=" "21×95="

Replace-Text "80×54=" "55×97="
Replace-Text "76×89=" "66×92="
Replace-Text "71×55=" "17×99="
Replace-Text "99×92=" "46×48="
Replace-Text "18×13=" "52×92="

Replace-Text "93×44=" "27×87="
Replace-Text "73×42=" "11×64="
Replace-Text "55×91=" "69×70="
Replace-Text "79×86=" "75×75="
Replace-Text "95×19=" "41×98="

Replace-Text "67×80=" "64×28="
Replace-Text "62×76=" "77×68="
Replace-Text "71×36=" "85×25="
Replace-Text "58×72=" "21×82="
Replace-Text "75×90=" "61×36="
